# The "2024" sheet keeps a stacked activity log per month: new entries are always
# inserted at the top of their month's block (columns R:S for September, rows 31-75;
# columns P:Q for August, rows 76-79), pushing older entries down by one row each time
# a new transaction is recorded. A brand-new September entry was logged, so the whole
# September block shifts down into row 76 (displacing the oldest August row into the
# August block), the August block shifts down into row 80 (displacing the "Broadband"
# group label that lived there), and "Broadband" is relocated to the newly added row 81.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert the new September transaction at the top of the log, which cascades everything
# else down by one row.
$ws.Range("R31").Value = "coimbatore ramalinga"
$ws.Range("S31").Value = "2024-09-05 17:06:01"

# Row 32
$ws.Range("R32").Value = "beneficiary"
$ws.Range("S32").Value = "2024-09-05 17:04:10"

# Row 33
$ws.Range("R33").Value = "bal axisbank"
$ws.Range("S33").Value = "2024-09-05 16:52:25"

# Row 34
$ws.Range("R34").Value = "share anyone axis"
$ws.Range("S34").Value = "2024-09-05 16:38:59"

# Row 35
$ws.Range("R35").Value = "transfer anyone axis"
$ws.Range("S35").Value = "2024-09-05 16:35:58"

# Row 36
$ws.Range("R36").Value = "share anyone axis"
$ws.Range("S36").Value = "2024-09-05 16:31:34"

# Row 37
$ws.Range("R37").Value = "transfer"
$ws.Range("S37").Value = "2024-09-05 16:28:38"

# Row 38
$ws.Range("R38").Value = "bal axisbank axis"
$ws.Range("S38").Value = "2024-09-05 16:26:56"

# Row 39
$ws.Range("R39").Value = "bal axisbank"
$ws.Range("S39").Value = "2024-09-05 16:26:55"

# Row 40
$ws.Range("S40").Value = "2024-09-05 16:25:07"

# Row 41
$ws.Range("R41").Value = "transfer"
$ws.Range("S41").Value = "2024-09-05 16:22:23"

# Row 42
$ws.Range("R42").Value = "share anyone axis"
$ws.Range("S42").Value = "2024-09-05 16:06:05"

# Row 43
$ws.Range("R43").Value = "internet bal axisbank"
$ws.Range("S43").Value = "2024-09-05 16:05:55"

# Row 44
$ws.Range("R44").Value = "transfer share anyone axis"
$ws.Range("S44").Value = "2024-09-05 16:03:14"

# Row 45
$ws.Range("R45").Value = "axis"

# Row 46
$ws.Range("R46").Value = "your net internet"
$ws.Range("S46").Value = "2024-09-05 15:57:15"

# Row 47
$ws.Range("R47").Value = "hear your feedback atm"
$ws.Range("S47").Value = "2024-09-05 14:21:08"

# Row 48
$ws.Range("S48").Value = "2024-09-05 14:18:32"

# Row 49
$ws.Range("S49").Value = "2024-09-05 14:13:16"

# Row 50
$ws.Range("R50").Value = "axis bna"
$ws.Range("S50").Value = "2024-09-05 14:15:23"

# Row 51
$ws.Range("R51").Value = "balance your axis"
$ws.Range("S51").Value = "2024-09-05 09:20:57"

# Row 52
$ws.Range("R52").Value = "bal axis"
$ws.Range("S52").Value = "2024-09-05 09:06:25"

# Row 53
$ws.Range("R53").Value = "broker"
$ws.Range("S53").Value = "2024-09-04 21:20:47"

# Row 54
$ws.Range("R54").Value = "exclusive on axis"
$ws.Range("S54").Value = "2024-09-04 13:21:05"

# Row 55
$ws.Range("R55").Value = "your corporate axis"
$ws.Range("S55").Value = "2024-09-04 11:46:10"

# Row 56
$ws.Range("R56").Value = "balance your axis"
$ws.Range("S56").Value = "2024-09-04 08:14:16"

# Row 57
$ws.Range("R57").Value = "axis"
$ws.Range("S57").Value = "2024-09-04 07:02:13"

# Row 58
$ws.Range("R58").Value = "bal axisbank w axis"
$ws.Range("S58").Value = "2024-09-04 06:53:15"

# Row 59
$ws.Range("R59").Value = "logging iob internet"
$ws.Range("S59").Value = "2024-09-03 20:09:12"

# Row 60
$ws.Range("R60").Value = "password internet"
$ws.Range("S60").Value = "2024-09-03 20:05:31"

# Row 61
$ws.Range("R61").Value = "logging iob internet"
$ws.Range("S61").Value = "2024-09-03 20:05:09"

# Row 62
$ws.Range("R62").Value = "internet"
$ws.Range("S62").Value = "2024-09-03 19:58:18"

# Row 63
$ws.Range("S63").Value = "2024-09-03 19:54:49"

# Row 64
$ws.Range("R64").Value = "login internet invalid"
$ws.Range("S64").Value = "2024-09-03 19:56:17"

# Row 65
$ws.Range("R65").Value = "corporate internet share"
$ws.Range("S65").Value = "2024-09-03 19:22:58"

# Row 66
$ws.Range("R66").Value = "login sbi internet personal do not share anyone"
$ws.Range("S66").Value = "2024-09-03 19:17:10"

# Row 67
$ws.Range("R67").Value = "login internet personal share"
$ws.Range("S67").Value = "2024-09-03 19:13:40"

# Row 68
$ws.Range("R68").Value = "internet verify it"
$ws.Range("S68").Value = "2024-09-03 19:05:49"

# Row 69
$ws.Range("R69").Value = "balance your axis"
$ws.Range("S69").Value = "2024-09-03 13:14:06"

# Row 70
$ws.Range("R70").Value = "lounge"
$ws.Range("S70").Value = "2024-09-03 13:08:08"

# Row 71
$ws.Range("R71").Value = "balance your axis"
$ws.Range("S71").Value = "2024-09-03 11:21:30"

# Row 72
$ws.Range("R72").Value = "broker"
$ws.Range("S72").Value = "2024-09-01 22:35:38"

# Row 73
$ws.Range("S73").Value = "2024-09-01 10:12:03"

# Row 74
$ws.Range("S74").Value = "2024-09-01 09:42:38"

# Row 75
$ws.Range("S75").Value = "2024-09-01 09:29:24"

# Row 76
$ws.Range("P76").Value = $null
$ws.Range("Q76").Value = $null
$ws.Range("R76").Value = "amazeloan"
$ws.Range("S76").Value = "2024-09-01 09:27:06"

# Row 77
$ws.Range("Q77").Value = "2024-08-30 12:15:48"

# Row 78
$ws.Range("Q78").Value = "2024-08-21 20:17:10"

# Row 79
$ws.Range("Q79").Value = "2024-08-21 20:16:45"

# Row 80
$ws.Range("A80").Value = $null
$ws.Range("P80").Value = "hdfc"
$ws.Range("Q80").Value = "2024-08-21 20:15:50"

# Row 81
$ws.Range("A81").Value = "Broadband"
